$ws = $excel.ActiveWorkbook.ActiveSheet

$values = @(
    0.093312672893958393
    -0.022450511379972227
    -0.0089999998327634501
    0.2839963621714432
    -0.0059999998388757803
    -0.0059999998332251891
    -0.026959134408373941
    -0.019999999795844658
    -0.0059999998266171417
    -0.0059999998237501018
    -0.0044999998269972252
    -0.0059999998226047957
    -0.0059999998195792159
    -0.011999999804113592
    -0.0059999998179245395
    -0.0059999998172621805
    -0.0059999998163950963
    -0.008999999808903425
    -0.0089999998323442298
    -0.0089999998280596571
    -0.0089999998271217407
    -0.050315594525227159
    -0.0089999998311762752
    -0.041999999748457029
    -0.041999999747050154
    -0.0059999998329303139
    -0.0059999998326185633
    -0.0059999998313937652
    -0.01199999981615818
    0.03183344881399286
    -0.014999999807018938
    -0.020999999792212343
    -0.0059999998284183675
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
